# Fruta / hortaliza, semanal
# Insert a new weekly record as row 42 in the data table (Arandano blue, Talca),
# pushing the existing rows 42-83 down to 43-84.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 42; Excel shifts rows 42:83 down to 43:84
# and copies formatting (including the date number format on column D)
# from the row immediately above.
$ws.Rows.Item(42).Insert()

# Populate the newly inserted row 42 with the new weekly observation.
$ws.Cells.Item(42, 1).Value = 5
$ws.Cells.Item(42, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(42, 3).Value = "Maule"
$ws.Cells.Item(42, 4).Value = 44893
$ws.Cells.Item(42, 5).Value = 7
$ws.Cells.Item(42, 6).Value = "Fruta"
$ws.Cells.Item(42, 7).Value = 100101
$ws.Cells.Item(42, 8).Value = "Berries"
$ws.Cells.Item(42, 9).Value = 100101001
$ws.Cells.Item(42, 10).Value = "Arándano (blue)"
$ws.Cells.Item(42, 11).Value = "Sin especificar"
$ws.Cells.Item(42, 12).Value = "Primera"
$ws.Cells.Item(42, 13).Value = 100
$ws.Cells.Item(42, 14).Value = 4000
$ws.Cells.Item(42, 15).Value = 4000
$ws.Cells.Item(42, 16).Value = 4000
$ws.Cells.Item(42, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(42, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(42, 19).Value = 2000
$ws.Cells.Item(42, 20).Value = 2
